# [Improvement] On terminology : room -> bed
$wb = $excel.ActiveWorkbook

# Rename the "rooms" sheet to "beds"
$roomsSheet = $wb.Worksheets.Item("rooms")
$roomsSheet.Name = "beds"

# Update the header row terminology on the beds sheet
$roomsSheet.Range("A1").Value = "all_beds"
$roomsSheet.Range("B1").Value = "new_beds"
$roomsSheet.Range("C1").Value = "old_beds"
$roomsSheet.Range("E1").Value = "new_beds_service"
$roomsSheet.Range("F1").Value = "old_beds_service"
$roomsSheet.Range("G1").Value = "beds_capacities"

# Make the beds sheet the active tab, with a new selection
$roomsSheet.Activate()
$roomsSheet.Range("E21").Select()
